$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 12
$ws.Range("B3").Value = "testwaiter@gmail.com"
$ws.Range("C3").Value = "Test"
$ws.Range("D3").Value = "Waiter"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = "Waiter"
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = $true
